$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 57,3
$data[0,0] = 43057
$data[0,1] = 6265
$data[0,2] = 7147
$data[1,0] = 22227
$data[1,1] = 2517
$data[1,2] = 2460
$data[2,0] = 74547
$data[2,1] = 6400
$data[2,2] = 6628
$data[3,0] = 1684
$data[3,1] = 636
$data[3,2] = 157
$data[4,0] = 47410
$data[4,1] = 7599
$data[4,2] = 6165
$data[5,0] = 5346
$data[5,1] = 1304
$data[5,2] = 1088
$data[6,0] = 5851
$data[6,1] = 1131
$data[6,2] = 681
$data[7,0] = 2795
$data[7,1] = 383
$data[7,2] = 266
$data[8,0] = 356
$data[8,1] = 199
$data[8,2] = 29
$data[9,0] = 3
$data[9,1] = 0
$data[9,2] = 0
$data[10,0] = 903
$data[10,1] = 225
$data[10,2] = 275
$data[11,0] = 2992
$data[11,1] = 1159
$data[11,2] = 854
$data[12,0] = 5030
$data[12,1] = 1885
$data[12,2] = 803
$data[13,0] = 3332
$data[13,1] = 1314
$data[13,2] = 522
$data[14,0] = 2192
$data[14,1] = 863
$data[14,2] = 182
$data[15,0] = 17078
$data[15,1] = 2657
$data[15,2] = 3100
$data[16,0] = 1457
$data[16,1] = 591
$data[16,2] = 448
$data[17,0] = 18668
$data[17,1] = 2200
$data[17,2] = 3061
$data[18,0] = 151
$data[18,1] = 363
$data[18,2] = 16
$data[19,0] = 16816
$data[19,1] = 2249
$data[19,2] = 2950
$data[20,0] = 1101
$data[20,1] = 407
$data[20,2] = 147
$data[21,0] = 19499
$data[21,1] = 2697
$data[21,2] = 3615
$data[22,0] = 76995
$data[22,1] = 7272
$data[22,2] = 9161
$data[23,0] = 6014
$data[23,1] = 1882
$data[23,2] = 970
$data[24,0] = 0
$data[24,1] = 0
$data[24,2] = 0
$data[25,0] = 5328
$data[25,1] = 1067
$data[25,2] = 1276
$data[26,0] = 1436
$data[26,1] = 443
$data[26,2] = 279
$data[27,0] = 14866
$data[27,1] = 2551
$data[27,2] = 2494
$data[28,0] = 481
$data[28,1] = 169
$data[28,2] = 204
$data[29,0] = 2629
$data[29,1] = 1704
$data[29,2] = 375
$data[30,0] = 16502
$data[30,1] = 3152
$data[30,2] = 2614
$data[31,0] = 10268
$data[31,1] = 2938
$data[31,2] = 2082
$data[32,0] = 5858
$data[32,1] = 600
$data[32,2] = 1387
$data[33,0] = 56782
$data[33,1] = 5349
$data[33,2] = 5394
$data[34,0] = 8435
$data[34,1] = 2729
$data[34,2] = 1240
$data[35,0] = 24560
$data[35,1] = 1933
$data[35,2] = 2802
$data[36,0] = 1056
$data[36,1] = 917
$data[36,2] = 182
$data[37,0] = 2087
$data[37,1] = 448
$data[37,2] = 732
$data[38,0] = 2228
$data[38,1] = 272
$data[38,2] = 93
$data[39,0] = 8333
$data[39,1] = 482
$data[39,2] = 265
$data[40,0] = 236
$data[40,1] = 91
$data[40,2] = 66
$data[41,0] = 536
$data[41,1] = 38
$data[41,2] = 43
$data[42,0] = 1045
$data[42,1] = 14
$data[42,2] = 7
$data[43,0] = 3417
$data[43,1] = 893
$data[43,2] = 431
$data[44,0] = 11818
$data[44,1] = 3402
$data[44,2] = 2023
$data[45,0] = 31291
$data[45,1] = 3334
$data[45,2] = 4184
$data[46,0] = 15164
$data[46,1] = 3479
$data[46,2] = 1202
$data[47,0] = 10946
$data[47,1] = 1123
$data[47,2] = 1531
$data[48,0] = 31694
$data[48,1] = 2899
$data[48,2] = 4397
$data[49,0] = 4842
$data[49,1] = 553
$data[49,2] = 1237
$data[50,0] = 14236
$data[50,1] = 2976
$data[50,2] = 2327
$data[51,0] = 1922
$data[51,1] = 1248
$data[51,2] = 717
$data[52,0] = 2006
$data[52,1] = 1265
$data[52,2] = 147
$data[53,0] = 3605
$data[53,1] = 901
$data[53,2] = 995
$data[54,0] = 12257
$data[54,1] = 4583
$data[54,2] = 2454
$data[55,0] = 12378
$data[55,1] = 897
$data[55,2] = 474
$data[56,0] = 669111
$data[56,1] = 99168
$data[56,2] = 91451

$ws.Range("B3:D59").Value = $data

Write-Output "Updated B3:D59 with statistics values"
